$wb = $excel.ActiveWorkbook

# Data to populate into each of the new sheets
$values = @(
    @(1, "selamlar", "selamlar"),
    @(2, "merhaba", "merhaba"),
    @(3, "deneme", "deneme")
)

function Fill-SheetData($ws) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 1
        $ws.Cells.Item($row, 1).Value = $values[$i][0]
        $ws.Cells.Item($row, 2).Value = $values[$i][1]
        $ws.Cells.Item($row, 3).Value = $values[$i][2]
    }
}

$sheet1 = $wb.Worksheets.Item("Sheet1")

# Add "Manipulated_Data" sheet right after Sheet1
$ws1 = $wb.Worksheets.Add($null, $sheet1)
$ws1.Name = "Manipulated_Data"
Fill-SheetData($ws1)

# Add "Manipulated_Data2" sheet right after "Manipulated_Data"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Manipulated_Data2"
Fill-SheetData($ws2)
